# Update Nalco PDF (2025-09-02 06:49:07 UTC)
# Insert a new top data row (row 2) with the latest price circular, pushing
# the existing rows down by one. Row 6 (old row 5, 07-08-2025 entry) becomes
# the new bottom row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture the existing data rows (2-5) before we overwrite anything ---
#        NOTE: use Value() (method call syntax) to force evaluation of the
#        COM property right now; plain ".Value" (no parens) yields the
#        property descriptor instead of the live cell value in this shim.
$oldA2 = $ws.Range("A2").Value()
$oldB2 = $ws.Range("B2").Value()
$oldC2 = $ws.Range("C2").Value()
$oldD2 = $ws.Range("D2").Value()
$oldE2 = $ws.Range("E2").Value()
$oldF2 = $ws.Range("F2").Value()

$oldA3 = $ws.Range("A3").Value()
$oldB3 = $ws.Range("B3").Value()
$oldC3 = $ws.Range("C3").Value()
$oldD3 = $ws.Range("D3").Value()
$oldE3 = $ws.Range("E3").Value()
$oldF3 = $ws.Range("F3").Value()

$oldA4 = $ws.Range("A4").Value()
$oldB4 = $ws.Range("B4").Value()
$oldC4 = $ws.Range("C4").Value()
$oldD4 = $ws.Range("D4").Value()
$oldE4 = $ws.Range("E4").Value()
$oldF4 = $ws.Range("F4").Value()

$oldA5 = $ws.Range("A5").Value()
$oldB5 = $ws.Range("B5").Value()
$oldC5 = $ws.Range("C5").Value()
$oldD5 = $ws.Range("D5").Value()
$oldE5 = $ws.Range("E5").Value()
$oldF5 = $ws.Range("F5").Value()

# --- 2. Make row 6 exist and carry the same look (border/alignment/number
#        format) as the other data rows, by copying formats down from row 5.
$ws.Range("A5:F5").Copy() | Out-Null
$ws.Range("A6:F6").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Remove all existing hyperlinks; we will recreate them below so the
#        relationship targets line up correctly with their (possibly moved)
#        cells.
$ws.Hyperlinks.Delete() | Out-Null

# --- 4. Shift the previous rows 2-5 down into rows 3-6. ---
$ws.Range("A6").Value = $oldA5
$ws.Range("B6").Value = $oldB5
$ws.Range("C6").Value = $oldC5
$ws.Range("D6").Value = $oldD5
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = $oldE5
$ws.Range("F6").Value = $oldF5

$ws.Range("A5").Value = $oldA4
$ws.Range("B5").Value = $oldB4
$ws.Range("C5").Value = $oldC4
$ws.Range("D5").Value = $oldD4
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = $oldE4
$ws.Range("F5").Value = $oldF4

$ws.Range("A4").Value = $oldA3
$ws.Range("B4").Value = $oldB3
$ws.Range("C4").Value = $oldC3
$ws.Range("D4").Value = $oldD3
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = $oldE3
$ws.Range("F4").Value = $oldF3

$ws.Range("A3").Value = $oldA2
$ws.Range("B3").Value = $oldB2
$ws.Range("C3").Value = $oldC2
$ws.Range("D3").Value = $oldD2
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = $oldE2
$ws.Range("F3").Value = $oldF2

# --- 5. Write the brand-new top row (row 2) with the latest circular. ---
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 272.05
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01-09-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

# --- 6. Restore each E-column cell's formatting to match the other data
#        cells (copy format from the header-adjacent data style cell A2,
#        which still carries the original "centered data" style).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("E2:E6").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = 0

# --- 7. Re-create the hyperlinks for F2:F6 pointing at their circular PDFs.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null

# --- 8. Adding hyperlinks applies Excel's default "Hyperlink" style (blue /
#        underlined font) to the cells. The source workbook keeps the plain
#        data-row style on column F, so restore it by copying the format
#        from the (still plainly styled) column E cells back onto column F.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2:F6").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = 0
